$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E35").Value = "V"
$ws.Range("F35").Value = "V"
$ws.Range("G35").Value = "V"
$ws.Range("H35").Value = "X"
$ws.Range("I35").Value = "V"
$ws.Range("J35").Value = "V"
$ws.Range("K35").Value = "X"
$ws.Range("L35").Value = "V"
